# Automatische test-sync: 2025-06-19 21:26:50
#
# This script applies the mail-log sync for the new "Afmelding nieuwsbrief"
# e-mail:
#   1. Appends a new row (17) to the "Logs" sheet with the incoming message.
#   2. Extends the conditional-formatting ranges on "Logs" to include row 17.
#   3. Re-orders the "Dashboard" summary rows 5/6 (Factuur / Administratie
#      now comes before Openingstijden / Locatie) and appends a new
#      "Afmelding / Nieuwsbrief" summary row (10).
#   4. Points the Dashboard bar chart's category/value series at the
#      extended A2:A10 / B2:B10 ranges.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- 1. New row on "Logs" ------------------------------------------------
$logs.Range("A17").Value = "Afmelding nieuwsbrief"
$logs.Range("B17").Value = "mailmind.test@zohomail.eu"
$logs.Range("C17").Value = "Graag afmelden voor de nieuwsbrief. Dank u."
$logs.Range("D17").Value = "Afmelding / Nieuwsbrief"
$logs.Range("F17").Value = "2025-06-19 21:26:27"
$logs.Range("G17").Value = "Nee"

# --- 2. Extend conditional formatting ranges on "Logs" --------------------
$catFormat = $logs.Range("D2:D16").FormatConditions.Item(1)
$catFormat.ModifyAppliesToRange($logs.Range("D2:D17"))

$answeredFormat = $logs.Range("G2:G16").FormatConditions.Item(1)
$answeredFormat.ModifyAppliesToRange($logs.Range("G2:G17"))

# --- 3. Swap Dashboard rows 5 & 6, append new summary row -----------------
$dashboard.Range("A5").Value = "Factuur / Administratie"
$dashboard.Range("A6").Value = "Openingstijden / Locatie"

$dashboard.Range("A10").Value = "Afmelding / Nieuwsbrief"
$dashboard.Range("B10").Value = 1

# --- 4. Update chart series ranges to include the new row -----------------
$chartObj = $dashboard.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES('Dashboard'!B1,'Dashboard'!`$A`$2:`$A`$10,'Dashboard'!`$B`$2:`$B`$10,1)"
